$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.756.75'
$ws.Range("E2").Value = '  -0.82%  '

$ws.Range("D3").Value = '2.231.44'
$ws.Range("E3").Value = '  -1.97%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''313.33'
$ws.Range("E5").Value = '  -1.56%  '

$ws.Range("D6").Value = '''98.16'
$ws.Range("E6").Value = '  -4.99%  '

$ws.Range("D7").Value = '''0.569'
$ws.Range("E7").Value = '  -3.34%  '

$ws.Range("E8").Value = '  +0.18%  '

$ws.Range("E9").Value = '  -6.87%  '

$ws.Range("D10").Value = '''35.71'
$ws.Range("E10").Value = '  -7.76%  '

$ws.Range("D11").Value = '''0.0820'
$ws.Range("E11").Value = '  -2.36%  '

$ws.Range("D12").Value = '''7.34'
$ws.Range("E12").Value = '  -6.63%  '

$ws.Range("E13").Value = '  -3.04%  '

$ws.Range("D14").Value = '2.570.71'
$ws.Range("E14").Value = '  -1.98%  '

$ws.Range("D15").Value = '2.233.79'
$ws.Range("E15").Value = '  -1.73%  '

$ws.Range("D16").Value = '''0.837'
$ws.Range("E16").Value = '  -4.32%  '

$ws.Range("D17").Value = '''14.00'

$ws.Range("D18").Value = '43.623.59'
$ws.Range("E18").Value = '  -0.87%  '

$ws.Range("D19").Value = '''13.03'
$ws.Range("E19").Value = '  -9.14%  '

$ws.Range("D20").Value = '0.0₃0965'
$ws.Range("E20").Value = '  -3.49%  '

$ws.Range("D21").Value = '''6.29'
$ws.Range("E21").Value = '  -5.60%  '

$ws.Range("D22").Value = '''65.14'
$ws.Range("E22").Value = '  -1.57%  '

$ws.Range("D23").Value = '''235.52'
$ws.Range("E23").Value = '  -1.01%  '

$ws.Range("D24").Value = '''2.97'
$ws.Range("E24").Value = '  -7.48%  '

$ws.Range("D25").Value = '''2.02'
$ws.Range("E25").Value = '  -8.44%  '

$ws.Range("E26").Value = '  +0.48%  '

$ws.Range("D27").Value = '''9.97'
$ws.Range("E27").Value = '  -2.81%  '

$ws.Range("D28").Value = '''2.17'
$ws.Range("E28").Value = '  -2.18%  '

$ws.Range("D29").Value = '''36.54'
$ws.Range("E29").Value = '  -7.17%  '

$ws.Range("D30").Value = '''5.97'
$ws.Range("E30").Value = '  -8.45%  '

$ws.Range("D31").Value = '''157.13'
$ws.Range("E31").Value = '  -3.04%  '

$ws.Range("D32").Value = '''19.85'
$ws.Range("E32").Value = '  -3.14%  '

$ws.Range("D33").Value = '''0.0825'
$ws.Range("E33").Value = '  -6.12%  '

$ws.Range("D34").Value = '''2.63'
$ws.Range("E34").Value = '  -3.28%  '

$ws.Range("D35").Value = '''3.10'
$ws.Range("E35").Value = '  -5.06%  '

$ws.Range("E36").Value = '  +0.08%  '

$ws.Range("E37").Value = '  -7.79%  '

$ws.Range("E38").Value = '  -3.41%  '

$ws.Range("D39").Value = '''15.45'
$ws.Range("E39").Value = '  -0.97%  '

$ws.Range("D40").Value = '''3.52'
$ws.Range("E40").Value = '  -8.60%  '

$ws.Range("D41").Value = '''4.00'
$ws.Range("E41").Value = '  -11.49%  '

$ws.Range("D42").Value = '''0.0306'
$ws.Range("E42").Value = '  -6.30%  '

$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("D44").Value = '1.710.04'
$ws.Range("E44").Value = '  -3.88%  '

$ws.Range("D45").Value = '''83.46'
$ws.Range("E45").Value = '  -1.75%  '

$ws.Range("E46").Value = '  -7.00%  '

$ws.Range("D47").Value = '''5.12'
$ws.Range("E47").Value = '  -5.52%  '

$ws.Range("D48").Value = '''101.43'
$ws.Range("E48").Value = '  -2.90%  '

$ws.Range("D49").Value = '''71.16'
$ws.Range("E49").Value = '  -4.71%  '

$ws.Range("E50").Value = '  +0.46%  '

$ws.Range("D51").Value = '''55.86'
$ws.Range("E51").Value = '  -6.20%  '
